$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# A fresh handoff-xliff generation pass just completed for the four
# "Ready for handoff" source files: their priority flips from "low" to
# "ht" (high) on every locale sheet, each locale's "Latest Handoff
# Datetime" is stamped with the new generation time, and the Overview
# rollup ("Latest HO Xliff Generate Date") reflects the same refresh.

$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4").Value = "2016-08-13 20:39:35"
$zhcn.Range("H5").Value = "2016-08-13 20:39:35"
$zhcn.Range("H6").Value = "2016-08-13 20:39:35"
$zhcn.Range("H7").Value = "2016-08-13 20:39:35"

$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4").Value = "2016-08-13 20:39:43"
$dede.Range("H5").Value = "2016-08-13 20:39:43"
$dede.Range("H6").Value = "2016-08-13 20:39:43"
$dede.Range("H7").Value = "2016-08-13 20:39:43"

$overview.Range("G4").Value = "2016-08-13 20:39:43"
$overview.Range("G5").Value = "2016-08-13 20:39:43"
$overview.Range("G6").Value = "2016-08-13 20:39:43"
$overview.Range("G7").Value = "2016-08-13 20:39:43"
